$d = $word.ActiveDocument

# --- 1) Split the "postalcode" run into "posta" + "lcode" (no text change, just a run split) ---
$rng = $d.Content
$found = $rng.Find.Execute("postalcode", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $rng.Start + 5
$headRng = $d.Range($rng.Start, $splitPoint)
$headRng.Bold = 1
$headRng.Bold = 0

# --- 2) Replace the trailing sentence with a period ---
$d.Content.Find.Execute(" and now our data is ready!", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)

# --- 3) Append the new paragraphs (exact OOXML) at the very end of the document ---
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newParagraphsXml = '<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:eastAsia="Leelawadee UI" w:hAnsi="Leelawadee UI" w:cs="Leelawadee UI"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:eastAsia="Leelawadee UI" w:hAnsi="Leelawadee UI" w:cs="Leelawadee UI"/></w:rPr><w:lastRenderedPageBreak/><w:t>Additionaly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:eastAsia="Leelawadee UI" w:hAnsi="Leelawadee UI" w:cs="Leelawadee UI"/></w:rPr><w:t xml:space="preserve"> we can get our venue data using Foursquare API.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve">First the data related to shopping venues in each neighborhood </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve">is </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve"> extracted</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t>using foursquare to check the neighborhoods with large scope for shopping.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve">Next, we analyzed each </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t>neighbourhood</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve"> by grouping the rows by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t>neighbourhood</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve"> and taking the mean of the frequency of occurrence of each venue category. Since we are analyzing the “Restaurants” data, we will filter the “Restaurants” as venue category for the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t>neighbourhoods</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve">On this we will run </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t>kmean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t xml:space="preserve"> cluster to get restaurant concentration in each neighborhood</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr><w:t>..</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:hAnsi="Leelawadee UI"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Leelawadee UI" w:eastAsia="Leelawadee UI" w:hAnsi="Leelawadee UI" w:cs="Leelawadee UI"/></w:rPr></w:pPr></w:p>'
$endRng.InsertXML($newParagraphsXml)
